$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10")
$ws.Range("A38").Value = "Nguyễn Minh Thảo"
$ws.Range("B38").Value = "Danh mục các món ăn nên gom nhóm lại theo dạng giống menu: món chính, món phụ, món uống…"
$ws.Range("B39").Value = "Thiếu thông tin thời gian phục vụ của nhà hàng"
$ws.Range("A41").Value = "Nguyễn Chí Hiếu"
$ws.Range("B41").Value = "Trang web thiết kế khó sử dụng, cần hoàn thiện các chức năng hơn, trang web cũng không thấy có gì đặc sắc"
$ws.Range("B41").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 25
